$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new date header in column CL (column 90), row 1
$ws.Cells.Item(1, 90).Value = "22-nov"

# Add the new counts for column CL, rows 2-11
$ws.Cells.Item(2, 90).Value = 13
$ws.Cells.Item(3, 90).Value = 8
$ws.Cells.Item(4, 90).Value = 8
$ws.Cells.Item(5, 90).Value = 12
$ws.Cells.Item(6, 90).Value = 9
$ws.Cells.Item(7, 90).Value = 12
$ws.Cells.Item(8, 90).Value = 10
$ws.Cells.Item(9, 90).Value = 11
$ws.Cells.Item(10, 90).Value = 17
$ws.Cells.Item(11, 90).Value = 0

# Match formatting of the previous column (CK) for the new CL column
$ws.Range("CK1").Copy()
$ws.Range("CL1").PasteSpecial(-4122)
$ws.Range("CK2:CK11").Copy()
$ws.Range("CL2:CL11").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Update the selection to match the new active cell
$ws.Range("CL11").Select()
